# Apply cryptos.xlsx price/volume updates (Tue Feb  6 19:52:46 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.181.99"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.373.42"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'303.47"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'97.21"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "'34.19"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("D13").Value = "'18.59"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "2.741.41"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "2.413.20"
$ws.Range("E16").Value = "  +4.70%  "
$ws.Range("D17").Value = "'0.804"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "43.167.86"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "'12.31"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  +4.97%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'68.27"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'236.01"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'24.82"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "'31.57"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.12"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "'17.28"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D36").Value = "'4.38"
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'22.71"
$ws.Range("E39").Value = "  +11.56%  "
$ws.Range("D40").Value = "'2.79"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "1.946.29"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0280"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'101.67"
$ws.Range("E44").Value = "  -38.61%  "
$ws.Range("E45").Value = "  +4.00%  "
$ws.Range("D46").Value = "'9.43"
$ws.Range("E46").Value = "  -10.07%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").Value = "2.597.66"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D49").Value = "'53.06"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "'1.51"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("D51").Value = "'72.34"
$ws.Range("E51").Value = "  +1.14%  "

Write-Output "Applied cryptos update."
